$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest data row (pulled 7th October 2024): CVR 42443611,
# Visma Løn, year 2024, TCV 83544, terminated 2024-07-11, quarter 2024Q3,
# same TCV_range bucket "80000-100000".

# Column A stores CVR numbers as text (as all the other rows in this sheet
# do), but simply assigning a digit-string to .Value lets Excel coerce it
# to a number. Stage the text value on an out-of-the-way scratch cell
# (forced to Text via NumberFormat), copy it, and paste values-only into
# A18 so the destination cell picks up the text value but none of the
# scratch cell's number-format/style baggage. Then wipe the scratch cell.
$scratch = $ws.Cells.Item(100, 100)
$scratch.NumberFormat = "@"
$scratch.Value = "42443611"
$scratch.Copy()
$ws.Cells.Item(18, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item(18, 2).Value = 2024
$ws.Cells.Item(18, 3).Value = 83544
$ws.Cells.Item(18, 4).Value = "Visma Løn"

# Column E uses the same custom date/time display format as the rows above.
$ws.Cells.Item(18, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 45484

$ws.Cells.Item(18, 8).Value = "2024Q3"
$ws.Cells.Item(18, 9).Value = "80000-100000"
